$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.622.93"
$ws.Range("E2").Value = "'  +0.08%  "
$ws.Range("D3").Value = "'2.466.47"
$ws.Range("E3").Value = "'  -0.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.27%  "
$ws.Range("D5").Value = "'317.71"
$ws.Range("E5").Value = "'  +0.99%  "
$ws.Range("D6").Value = "'92.34"
$ws.Range("E6").Value = "'  +1.23%  "
$ws.Range("D7").Value = "'0.551"
$ws.Range("E7").Value = "'  +0.31%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "'  +0.44%  "
$ws.Range("D10").Value = "'0.0869"
$ws.Range("E10").Value = "'  +9.64%  "
$ws.Range("D11").Value = "'32.92"
$ws.Range("E11").Value = "'  +1.08%  "
$ws.Range("D12").Value = "'0.110"
$ws.Range("E12").Value = "'  +0.21%  "
$ws.Range("D13").Value = "'2.847.31"
$ws.Range("E13").Value = "'  -0.01%  "
$ws.Range("D14").Value = "'6.87"
$ws.Range("E14").Value = "'  +0.07%  "
$ws.Range("D15").Value = "'15.49"
$ws.Range("E15").Value = "'  -2.11%  "
$ws.Range("D16").Value = "'2.491.84"
$ws.Range("E16").Value = "'  +0.75%  "
$ws.Range("D17").Value = "'0.795"
$ws.Range("E17").Value = "'  +2.32%  "
$ws.Range("D18").Value = "'41.595.75"
$ws.Range("E18").Value = "'  -0.02%  "
$ws.Range("D19").Value = "'6.45"
$ws.Range("E19").Value = "'  -0.79%  "
$ws.Range("D20").Value = "'0.0₃0946"
$ws.Range("E20").Value = "'  +0.47%  "
$ws.Range("D21").Value = "'70.83"
$ws.Range("E21").Value = "'  -0.55%  "
$ws.Range("D22").Value = "'11.26"
$ws.Range("E22").Value = "'  +0.51%  "
$ws.Range("D23").Value = "'239.87"
$ws.Range("E23").Value = "'  +0.69%  "
$ws.Range("D24").Value = "'2.74"
$ws.Range("E24").Value = "'  +0.64%  "
$ws.Range("D25").Value = "'1.95"
$ws.Range("E25").Value = "'  +2.03%  "
$ws.Range("E26").Value = "'  +0.01%  "
$ws.Range("D27").Value = "'24.84"
$ws.Range("E27").Value = "'  +0.80%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "'  -0.07%  "
$ws.Range("D29").Value = "'9.71"
$ws.Range("E29").Value = "'  +0.36%  "
$ws.Range("D30").Value = "'36.75"
$ws.Range("E30").Value = "'  +4.13%  "
$ws.Range("D31").Value = "'157.68"
$ws.Range("E31").Value = "'  +1.25%  "
$ws.Range("D32").Value = "'5.48"
$ws.Range("E32").Value = "'  +0.57%  "
$ws.Range("E33").Value = "'  +0.05%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "'  -0.48%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0760"
$ws.Range("E35").Value = "'  +0.26%  "
$ws.Range("D36").Value = "'17.23"
$ws.Range("E36").Value = "'  -0.07%  "
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "'  +3.39%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.116"
$ws.Range("E38").Value = "'  +1.47%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.90"
$ws.Range("E39").Value = "'  +0.58%  "
$ws.Range("E40").Value = "'  +1.20%  "
$ws.Range("D41").Value = "'2.54"
$ws.Range("E41").Value = "'  +8.18%  "
$ws.Range("D42").Value = "'3.99"
$ws.Range("E42").Value = "'  -0.24%  "
$ws.Range("D43").Value = "'1.987.05"
$ws.Range("E43").Value = "'  +1.33%  "
$ws.Range("D44").Value = "'18.97"
$ws.Range("E44").Value = "'  +2.16%  "
$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "'  +0.01%  "
$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "'  +1.51%  "
$ws.Range("D47").Value = "'9.41"
$ws.Range("E47").Value = "'  +3.92%  "
$ws.Range("D48").Value = "'2.710.61"
$ws.Range("E48").Value = "'  +0.21%  "
$ws.Range("D49").Value = "'97.55"
$ws.Range("E49").Value = "'  +0.48%  "
$ws.Range("D50").Value = "'75.77"
$ws.Range("E50").Value = "'  +5.95%  "
$ws.Range("D51").Value = "'66.56"
$ws.Range("E51").Value = "'  -1.15%  "
